$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update text/link/percentage cells and non-ambiguous price strings ---
$ws.Range("D2").Value = '27.177.94'
$ws.Range("E2").Value = '  -3.57%  '
$ws.Range("D3").Value = '1.805.25'
$ws.Range("E3").Value = '  -3.80%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  -2.44%  '
$ws.Range("E8").Value = '  -3.78%  '
$ws.Range("E9").Value = '  -4.07%  '
$ws.Range("E10").Value = '  -4.06%  '
$ws.Range("E11").Value = '  -4.65%  '
$ws.Range("D12").Value = '1.803.01'
$ws.Range("E12").Value = '  -5.04%  '
$ws.Range("E13").Value = '  -2.91%  '
$ws.Range("E14").Value = '  -3.87%  '
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("E18").Value = '  -4.11%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  -3.39%  '
$ws.Range("D21").Value = '27.424.66'
$ws.Range("E21").Value = '  -3.11%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").Value = '2.079.52'
$ws.Range("E24").Value = '  -4.15%  '
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("E28").Value = '  -6.88%  '
$ws.Range("E29").Value = '  -4.22%  '
$ws.Range("E30").Value = '  -8.92%  '
$ws.Range("E32").Value = '  -6.31%  '
$ws.Range("E33").Value = '  -5.45%  '
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("E35").Value = '  -6.23%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -5.87%  '
$ws.Range("E38").Value = '  -4.98%  '
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E40").Value = '  -3.44%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E41").Value = '  -6.54%  '
$ws.Range("E42").Value = '  -3.95%  '
$ws.Range("E43").Value = '  -8.81%  '
$ws.Range("E44").Value = '  -4.60%  '
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("E46").Value = '  -3.47%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("E49").Value = '  -4.58%  '
$ws.Range("E50").Value = '  -3.93%  '
$ws.Range("E51").Value = '  -4.41%  '

# --- Update purely numeric-looking price strings, forcing them to stay as text ---
# (Excel would otherwise silently convert these into numeric values)
$numericPriceCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D5").Value = '310.16'
$ws.Range("D6").Value = '1.003'
$ws.Range("D7").Value = '0.4204'
$ws.Range("D8").Value = '0.3544'
$ws.Range("D9").Value = '0.07107'
$ws.Range("D10").Value = '0.8464'
$ws.Range("D11").Value = '20.10'
$ws.Range("D13").Value = '5.308'
$ws.Range("D14").Value = '6.355'
$ws.Range("D15").Value = '0.06885'
$ws.Range("D16").Value = '1.006'
$ws.Range("D17").Value = '81.21'
$ws.Range("D18").Value = '0.000008753'
$ws.Range("D19").Value = '1.001'
$ws.Range("D20").Value = '15.04'
$ws.Range("D22").Value = '5.076'
$ws.Range("D23").Value = '10.84'
$ws.Range("D25").Value = '1.957'
$ws.Range("D26").Value = '153.40'
$ws.Range("D27").Value = '18.18'
$ws.Range("D28").Value = '5.031'
$ws.Range("D29").Value = '113.02'
$ws.Range("D30").Value = '1.707'
$ws.Range("D31").Value = '0.08864'
$ws.Range("D32").Value = '0.7373'
$ws.Range("D33").Value = '4.450'
$ws.Range("D34").Value = '2.913'
$ws.Range("D35").Value = '1.099'
$ws.Range("D36").Value = '1.003'
$ws.Range("D37").Value = '1.065'
$ws.Range("D38").Value = '0.05201'
$ws.Range("D39").Value = '0.01892'
$ws.Range("D40").Value = '0.1632'
$ws.Range("D41").Value = '2.705'
$ws.Range("D42").Value = '0.4953'
$ws.Range("D43").Value = '6.260'
$ws.Range("D44").Value = '8.153'
$ws.Range("D45").Value = '104.70'
$ws.Range("D46").Value = '10.21'
$ws.Range("D47").Value = '1.003'
$ws.Range("D48").Value = '0.06382'
$ws.Range("D49").Value = '0.4544'
$ws.Range("D50").Value = '1.587'
$ws.Range("D51").Value = '62.62'

foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
